$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow column B to make room for the new column C, and size column C
$ws.Columns.Item(2).ColumnWidth = 37.33
$ws.Columns.Item(3).ColumnWidth = 24

# New data block (rows 4-5) describing the transmitter actions / byte codes
$ws.Range("A4").Value = "actions"
$ws.Range("B4").Value = "0x01"
$ws.Range("C4").Value = "power "

$ws.Range("B5").Value = "0x02"
$ws.Range("C5").Value = "reset"

# Match the final selection from the diff
$ws.Range("C5").Select()
